$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for NATMI LR-pair output (Omg-Tnfrsf1b)
$newValues = @{
  "G2" = 3.070536666666667
  "H2" = 9.21161
  "I2" = 0.2245190988242715
  "J2" = 0.2245190988242715
  "M2" = 12.57753066666667
  "N2" = 37.732592
  "O2" = 0.1317204310459389
  "P2" = 0.1317204310459389
  "Q2" = 38.61976908812444
  "R2" = 347.57792179312
  "S2" = 0.02957375247517879
  "T2" = 0.02957375247517879
  "G3" = 3.070536666666667
  "H3" = 9.21161
  "I3" = 0.2245190988242715
  "J3" = 0.2245190988242715
  "O3" = 0.1145776761962127
  "P3" = 0.1145776761962127
  "Q3" = 33.59359943036
  "R3" = 302.34239487324
  "S3" = 0.02572487660495286
  "T3" = 0.02572487660495285
  "G4" = 3.070536666666667
  "H4" = 9.21161
  "I4" = 0.2245190988242715
  "J4" = 0.2245190988242715
  "M4" = 2.520244333333333
  "N4" = 7.560733
  "O4" = 0.02639370785296846
  "P4" = 0.02639370785296846
  "Q4" = 7.738502634458889
  "R4" = 69.64652371013
  "S4" = 0.005925891501779576
  "T4" = 0.005925891501779575
  "G5" = 3.070536666666667
  "H5" = 9.21161
  "I5" = 0.2245190988242715
  "J5" = 0.2245190988242715
  "M5" = 69.44815566666666
  "N5" = 208.344467
  "O5" = 0.7273081849048799
  "P5" = 0.7273081849048799
  "Q5" = 213.2431084068745
  "R5" = 1919.18797566187
  "S5" = 0.1632945782423603
  "T5" = 0.1632945782423602
  "I6" = 0.2851553493233187
  "J6" = 0.2851553493233187
  "M6" = 12.57753066666667
  "N6" = 37.732592
  "O6" = 0.1317204310459389
  "P6" = 0.1317204310459389
  "Q6" = 49.04987505641778
  "R6" = 441.44887550776
  "S6" = 0.03756078552792282
  "T6" = 0.03756078552792282
  "I7" = 0.2851553493233187
  "J7" = 0.2851553493233187
  "O7" = 0.1145776761962127
  "P7" = 0.1145776761962127
  "S7" = 0.03267243728038512
  "T7" = 0.03267243728038512
  "I8" = 0.2851553493233187
  "J8" = 0.2851553493233187
  "M8" = 2.520244333333333
  "N8" = 7.560733
  "O8" = 0.02639370785296846
  "P8" = 0.02639370785296846
  "Q8" = 9.828453051540556
  "R8" = 88.456077463865
  "S8" = 0.007526306982750841
  "T8" = 0.007526306982750841
  "I9" = 0.2851553493233187
  "J9" = 0.2851553493233187
  "M9" = 69.44815566666666
  "N9" = 208.344467
  "O9" = 0.7273081849048799
  "P9" = 0.7273081849048799
  "Q9" = 270.8340332157928
  "R9" = 2437.506298942135
  "S9" = 0.2073958195322599
  "T9" = 0.2073958195322599
  "G10" = 3.959514
  "H10" = 11.878542
  "I10" = 0.2895215434854775
  "J10" = 0.2895215434854775
  "M10" = 12.57753066666667
  "N10" = 37.732592
  "O10" = 0.1317204310459389
  "P10" = 0.1317204310459389
  "Q10" = 49.800908760096
  "R10" = 448.2081788408639
  "S10" = 0.03813590250499264
  "T10" = 0.03813590250499264
  "G11" = 3.959514
  "H11" = 11.878542
  "I11" = 0.2895215434854775
  "J11" = 0.2895215434854775
  "O11" = 0.1145776761962127
  "P11" = 0.1145776761962127
  "Q11" = 43.31956973479199
  "R11" = 389.8761276131279
  "S11" = 0.03317270566130675
  "T11" = 0.03317270566130675
  "G12" = 3.959514
  "H12" = 11.878542
  "I12" = 0.2895215434854775
  "J12" = 0.2895215434854775
  "M12" = 2.520244333333333
  "N12" = 7.560733
  "O12" = 0.02639370785296846
  "P12" = 0.02639370785296846
  "Q12" = 9.978942721254001
  "R12" = 89.81048449128599
  "S12" = 0.007641547035896196
  "T12" = 0.007641547035896196
  "G13" = 3.959514
  "H13" = 11.878542
  "I13" = 0.2895215434854775
  "J13" = 0.2895215434854775
  "M13" = 69.44815566666666
  "N13" = 208.344467
  "O13" = 0.7273081849048799
  "P13" = 0.7273081849048799
  "Q13" = 274.980944636346
  "R13" = 2474.828501727114
  "S13" = 0.2105713882832819
  "T13" = 0.2105713882832819
  "G14" = 2.746207666666667
  "H14" = 8.238623
  "I14" = 0.2008040083669322
  "J14" = 0.2008040083669322
  "M14" = 12.57753066666667
  "N14" = 37.732592
  "O14" = 0.1317204310459389
  "P14" = 0.1317204310459389
  "Q14" = 34.54051114453511
  "R14" = 310.864600300816
  "S14" = 0.02644999053784463
  "T14" = 0.02644999053784463
  "G15" = 2.746207666666667
  "H15" = 8.238623
  "I15" = 0.2008040083669322
  "J15" = 0.2008040083669322
  "O15" = 0.1145776761962127
  "P15" = 0.1145776761962127
  "Q15" = 30.045236491748
  "R15" = 270.407128425732
  "S15" = 0.02300765664956794
  "T15" = 0.02300765664956794
  "G16" = 2.746207666666667
  "H16" = 8.238623
  "I16" = 0.2008040083669322
  "J16" = 0.2008040083669322
  "M16" = 2.520244333333333
  "N16" = 7.560733
  "O16" = 0.02639370785296846
  "P16" = 0.02639370785296846
  "Q16" = 6.921114310073222
  "R16" = 62.290028790659
  "S16" = 0.005299962332541842
  "T16" = 0.005299962332541842
  "G17" = 2.746207666666667
  "H17" = 8.238623
  "I17" = 0.2008040083669322
  "J17" = 0.2008040083669322
  "M17" = 69.44815566666666
  "N17" = 208.344467
  "O17" = 0.7273081849048799
  "P17" = 0.7273081849048799
  "Q17" = 190.7190575276601
  "R17" = 1716.471517748941
  "S17" = 0.1460463988469778
  "T17" = 0.1460463988469778
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}
